# Update the header date.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2026-01-05 Monday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2026-01-06 Tuesday", 2)

# Update the 25 division problems in the table. We address cells by their
# (row, column) position rather than doing a global text Find/Replace,
# because several of the new values collide with *other* cells' old
# values (e.g. "23÷7=" becomes "38÷8=", while a different cell's original
# "38÷8=" becomes "36÷5=") -- a sequential Find/Replace would end up
# clobbering the freshly written text. Direct cell addressing sidesteps
# that ordering problem entirely.
$t = $d.Tables.Item(1)

$newValues = @(
    @{Row=1;  Col=1; Text="68÷2="},
    @{Row=1;  Col=2; Text="58÷4="},
    @{Row=1;  Col=3; Text="26÷4="},
    @{Row=1;  Col=4; Text="23÷4="},
    @{Row=1;  Col=5; Text="50÷5="},

    @{Row=5;  Col=1; Text="24÷8="},
    @{Row=5;  Col=2; Text="74÷8="},
    @{Row=5;  Col=3; Text="80÷4="},
    @{Row=5;  Col=4; Text="91÷6="},
    @{Row=5;  Col=5; Text="36÷5="},

    @{Row=9;  Col=1; Text="76÷5="},
    @{Row=9;  Col=2; Text="38÷8="},
    @{Row=9;  Col=3; Text="22÷3="},
    @{Row=9;  Col=4; Text="75÷5="},
    @{Row=9;  Col=5; Text="19÷6="},

    @{Row=13; Col=1; Text="64÷3="},
    @{Row=13; Col=2; Text="33÷7="},
    @{Row=13; Col=3; Text="36÷9="},
    @{Row=13; Col=4; Text="11÷4="},
    @{Row=13; Col=5; Text="12÷9="},

    @{Row=17; Col=1; Text="47÷8="},
    @{Row=17; Col=2; Text="77÷2="},
    @{Row=17; Col=3; Text="22÷4="},
    @{Row=17; Col=4; Text="62÷8="},
    @{Row=17; Col=5; Text="92÷9="}
)

foreach ($entry in $newValues) {
    $cell = $t.Cell($entry.Row, $entry.Col)
    $cell.Range.Text = $entry.Text
}

Write-Output "Updated date and $($newValues.Count) table cells."
